# Documentation Updates for Release (March 12, 2018 Release, v18.312)
#
# On the "Sample Interaction Mappings" sheet:
#   - Column D header "Interaction Source Field API Name" is renamed to
#     "Source Field API Name" (it now pairs with a new "Source Object API
#     Name" column).
#   - A new column G "Source Object API Name" is appended, populated with
#     "Interaction__c" for every mapping row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sample Interaction Mappings")

# Bring column G's formatting in line with the existing column F (header
# style on row 1, body style for the data rows) before filling in values.
$ws.Range("F1:F63").Copy()
$ws.Range("G1:G63").PasteSpecial(-4122)

# Rename the existing "Interaction Source Field API Name" header.
$ws.Cells.Item(1, 4).Value = "Source Field API Name"

# New column header.
$ws.Cells.Item(1, 7).Value = "Source Object API Name"

# Every data row (2-63) maps to the Interaction__c object.
$ws.Range("G2:G63").Value = "Interaction__c"

# Reflect the author's final on-screen state: the workbook was left open on
# the first sheet ("Sample Interaction Mappings") with D73 selected (rather
# than the "Sample Leads" sheet that was active before the edit).
[void]$ws.Activate()
[void]$ws.Range("D73").Select()
